# Natmi following Dr Hou advice
# Replace the Timp2-Itga3 LR-pair result table (rows 2-13) with the updated
# values, and append 4 new rows (14-17) completing the sending x target
# cluster grid (adds the previously-missing "M2" target-cluster rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 16 data rows (row 2..17) x 20 columns (A..T)
$rows = @(
  ,@("ECs", "Timp2", "Itga3", "ECs", 3, 1, 6.998221666666667, 20.994665, 0.01819620957294902, 0.01819620957294902, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 89.4297521444489, 804.8677693000401, 0.01311398013997092, 0.01311398013997092)
  ,@("ECs", "Timp2", "Itga3", "FAPs", 3, 1, 6.998221666666667, 20.994665, 0.01819620957294902, 0.01819620957294902, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 6.927033423132778, 62.34330080819501, 0.001015780280740907, 0.001015780280740907)
  ,@("ECs", "Timp2", "Itga3", "M2", 3, 1, 6.998221666666667, 20.994665, 0.01819620957294902, 0.01819620957294902, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 0.2082857387244445, 1.87457164852, 0.00003054302372055814, 0.00003054302372055814)
  ,@("ECs", "Timp2", "Itga3", "sCs", 3, 1, 6.998221666666667, 20.994665, 0.01819620957294902, 0.01819620957294902, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 27.52254318667167, 247.702888680045, 0.004035906128516628, 0.004035906128516628)
  ,@("FAPs", "Timp2", "Itga3", "ECs", 3, 1, 255.6993613333333, 767.098084, 0.6648487841778685, 0.6648487841778684, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 3267.563046259687, 29408.06741633718, 0.4791554920731408, 0.4791554920731408)
  ,@("FAPs", "Timp2", "Itga3", "FAPs", 3, 1, 255.6993613333333, 767.098084, 0.6648487841778685, 0.6648487841778684, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 253.0983021967302, 2277.884719770572, 0.0371143386722928, 0.0371143386722928)
  ,@("FAPs", "Timp2", "Itga3", "M2", 3, 1, 255.6993613333333, 767.098084, 0.6648487841778685, 0.6648487841778684, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 7.610294858243556, 68.492653724192, 0.001115973747407101, 0.001115973747407101)
  ,@("FAPs", "Timp2", "Itga3", "sCs", 3, 1, 255.6993613333333, 767.098084, 0.6648487841778685, 0.6648487841778684, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 1005.612146957481, 9050.509322617332, 0.1474629796850277, 0.1474629796850277)
  ,@("M2", "Timp2", "Itga3", "ECs", 3, 1, 84.71970766666666, 254.159123, 0.2202813270411758, 0.2202813270411758, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 1082.626818534161, 9743.641366807447, 0.1587564122320802, 0.1587564122320802)
  ,@("M2", "Timp2", "Itga3", "FAPs", 3, 1, 84.71970766666666, 254.159123, 0.2202813270411758, 0.2202813270411758, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 83.85791056037878, 754.721195043409, 0.01229692520999039, 0.01229692520999039)
  ,@("M2", "Timp2", "Itga3", "M2", 3, 1, 84.71970766666666, 254.159123, 0.2202813270411758, 0.2202813270411758, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 2.521484419380445, 22.693359774424, 0.0003697505115030534, 0.0003697505115030534)
  ,@("M2", "Timp2", "Itga3", "sCs", 3, 1, 84.71970766666666, 254.159123, 0.2202813270411758, 0.2202813270411758, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 333.1849038340977, 2998.664134506879, 0.04885823908760208, 0.04885823908760208)
  ,@("sCs", "Timp2", "Itga3", "ECs", 3, 1, 37.18048166666667, 111.541445, 0.09667367920800672, 0.09667367920800671, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 475.1265990757022, 4276.139391681319, 0.06967257131817341, 0.06967257131817341)
  ,@("sCs", "Timp2", "Itga3", "FAPs", 3, 1, 37.18048166666667, 111.541445, 0.09667367920800672, 0.09667367920800671, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 36.80226941365945, 331.220424722935, 0.005396685315833641, 0.005396685315833641)
  ,@("sCs", "Timp2", "Itga3", "M2", 3, 1, 37.18048166666667, 111.541445, 0.09667367920800672, 0.09667367920800671, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 1.106590282351111, 9.959312541160001, 0.0001622704149106609, 0.0001622704149106609)
  ,@("sCs", "Timp2", "Itga3", "sCs", 3, 1, 37.18048166666667, 111.541445, 0.09667367920800672, 0.09667367920800671, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 146.2230636743317, 1316.007573068985, 0.021442152159089, 0.021442152159089)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
  for ($j = 0; $j -lt $numCols; $j++) {
    $data[$i, $j] = $rows[$i][$j]
  }
}

$ws.Range("A2:T17").Value = $data
